# Algorithms tests comparison table in tries.xlsx
#
# Adds a new "10/11" result column (G) to the comparison table on Sheet1:
#   - G2 gets the new value "10/11"
#   - G3 gets the value "11/12" (already used elsewhere in the sheet)
# and leaves the G2 cell selected, matching the author's last saved
# selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "10/11"
$ws.Range("G3").Value = "11/12"

# Leave the selection on G2, as captured in the saved sheet view.
$null = $ws.Range("G2").Select()
